# "added buy many items from pega feature"
# The Shopping List sheet shrinks from 7 item rows down to 3 "bought in bulk"
# rows, with the kept rows re-pointed at different Products-sheet items and
# quantities, plus one brand-new line the customer is buying straight from
# Pega (order id 12345678, outside the normal product catalog run).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Shopping List")

# Drop the last four shopping-list rows (old rows 5-8) - delete bottom-up so
# earlier row numbers stay valid while we work.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()

# Row 2 -> Outback Lager, qty 5 (was Coffee Milk, qty 2)
$ws.Range("A2").Value2 = "Outback Lager"
$ws.Range("B2").Value2 = 5
$ws.Range("C2").Value2 = 127508
$ws.Range("D2").Value2 = 47414620

# Row 3 -> Veggie Spread, qty 3 (was Quohogs, qty 5)
$ws.Range("A3").Value2 = "Veggie Spread"
$ws.Range("B3").Value2 = 3
$ws.Range("C3").Value2 = 127508
$ws.Range("D3").Value2 = 37614806

# Row 4 -> Chef Anton's Cajun Seasoning, qty 4, brand-new Pega order id
$ws.Range("A4").Value2 = "Chef Anton's Cajun Seasoning"
$ws.Range("B4").Value2 = 4
$ws.Range("C4").Value2 = 127508
$ws.Range("D4").Value2 = 12345678

# Column A now needs to fit the longest item name ("Chef Anton's Cajun
# Seasoning"), so widen/re-bestfit it like Excel would after the edit.
$ws.Columns.Item(1).ColumnWidth = 22.21875
